# Update the dSF (column F) values for the rows whose source data was repulled.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -3
$ws.Range("F3").Value = -4
$ws.Range("F6").Value = -1
$ws.Range("F8").Value = -3
$ws.Range("F13").Value = -4
$ws.Range("F18").Value = -3
